$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '45.771.85'
$ws.Range('E2').Value = '  -3.84%  '

$ws.Range('D3').Value = '2.668.68'
$ws.Range('E3').Value = '  +0.86%  '

$ws.Range('E4').Value = '  +0.25%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '311.99'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.77%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '98.42'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -7.23%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.598'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -2.98%  '

$ws.Range('E8').Value = '  +0.23%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.582'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -3.38%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '38.36'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -5.34%  '

$ws.Range('E11').Value = '  -1.63%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '8.15'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -3.93%  '

$ws.Range('D13').Value = '3.078.49'
$ws.Range('E13').Value = '  +0.94%  '

$ws.Range('E14').Value = '  +0.51%  '

$ws.Range('D15').Value = '2.672.01'
$ws.Range('E15').Value = '  +0.49%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.933'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -2.08%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '15.13'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -2.16%  '

$ws.Range('D18').Value = '45.759.35'
$ws.Range('E18').Value = '  -4.22%  '

$ws.Range('E19').Value = '  -2.31%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.86'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.46%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.85'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -4.74%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '75.88'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +3.71%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '282.50'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +3.08%  '

$ws.Range('E24').Value = '  -1.95%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '31.16'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.25%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.24'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.64%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.997'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.25%  '

$ws.Range('E28').Value = '  -0.80%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '10.58'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -2.23%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '38.26'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -5.41%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.16'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -6.54%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.23'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.29%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.74'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.43%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.35'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +3.52%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '155.01'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.40%  '

$ws.Range('B36').Value = 'WEMIXToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.83'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.93%  '

$ws.Range('B37').Value = 'Hedera'
$ws.Range('C37').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0841'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.59%  '

$ws.Range('E38').Value = '  -2.37%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '25.75'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +10.49%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.124'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.00%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '16.30'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.24%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.61'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -4.35%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0328'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.25%  '

$ws.Range('E44').Value = '  -7.85%  '

$ws.Range('D45').Value = '2.121.50'
$ws.Range('E45').Value = '  -3.27%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.998'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.15%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '94.25'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.38%  '

$ws.Range('B48').Value = 'Aave'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '111.94'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -3.22%  '

$ws.Range('B49').Value = 'FraxShare'
$ws.Range('C49').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.31'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -7.79%  '

$ws.Range('D50').Value = '2.929.70'
$ws.Range('E50').Value = '  +0.96%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.200'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.76%  '
